$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 6 brand-new rows right after the header, pushing the
# existing data rows (old rows 2-21) down to rows 8-27. ---
$ws.Rows("2:7").Insert(-4121)

# The inserted rows pick up formatting from the row above (the bold header),
# so strip that back out to match the plain/unstyled data rows.
$ws.Range("A2:C7").ClearFormats()

# --- Step 2: populate the newly inserted rows (new rows 2-7) ---
$topRows = @(
    @(-1.428148408217626, 1.415199639612464, 1.873757413392585),
    @(-0.3362366334978593, 0.3509886704687749, 0.4294003237217614),
    @(1.29839596625184, 0.07410049284516276, -0.7659658990882879),
    @(1.272031672326387, -1.028642801671032, -0.5943050091996711),
    @(0.3833239707559676, -2.850368825271993, -0.3757410840134667),
    @(0.4730243273766702, -4.100493916726186, -1.476941543751531)
)

for ($i = 0; $i -lt $topRows.Count; $i++) {
    $r = 2 + $i
    $vals = $topRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}

# --- Step 3: append 4 brand-new rows at the bottom (new rows 28-31) ---
$bottomRows = @(
    @(2.691667430955088, 9.822626233541001, -7.046651205013665),
    @(4.03936266899109, -6.018842667231255, -2.086610792307836),
    @(3.831653899372287, 3.543843676683208, -1.997865703273085),
    @(-3.046866848899856, 0.7223788134726359, 2.788422576615722)
)

for ($i = 0; $i -lt $bottomRows.Count; $i++) {
    $r = 28 + $i
    $vals = $bottomRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
